# Auto-generated edit script: update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '66.255.96'),
    @('E2', '  +0.09%  '),
    @('D3', '3.566.65'),
    @('E3', '  +0.08%  '),
    @('E4', '  -0.06%  '),
    @('D5', '604.99'),
    @('E5', '  -0.09%  '),
    @('D6', '147.53'),
    @('E6', '  +2.32%  '),
    @('D7', '3.565.13'),
    @('E7', '  +0.07%  '),
    @('E8', '  -0.09%  '),
    @('D9', '0.489'),
    @('E9', '  +0.11%  '),
    @('D10', '0.135'),
    @('E10', '  -1.12%  '),
    @('D11', '7.90'),
    @('E11', '  +1.58%  '),
    @('E12', '  -0.70%  '),
    @('D13', '4.175.65'),
    @('E13', '  +0.13%  '),
    @('D14', '0.0000204'),
    @('E14', '  -1.16%  '),
    @('D15', '29.41'),
    @('E15', '  -3.03%  '),
    @('D16', '3.552.15'),
    @('E16', '  -0.49%  '),
    @('E17', '  +1.65%  '),
    @('D18', '66.281.15'),
    @('D19', '10.96'),
    @('E19', '  -3.91%  '),
    @('D20', '6.28'),
    @('E20', '  +1.00%  '),
    @('D21', '14.76'),
    @('E21', '  -0.14%  '),
    @('D22', '420.89'),
    @('E22', '  -2.18%  '),
    @('D23', '0.608'),
    @('E23', '  -0.76%  '),
    @('D24', '77.76'),
    @('E24', '  -2.19%  '),
    @('D25', '3.706.31'),
    @('E25', '  +0.05%  '),
    @('E26', '  +0.02%  '),
    @('E27', '  +0.89%  '),
    @('D28', '9.28'),
    @('E28', '  +1.46%  '),
    @('D29', '8.02'),
    @('E29', '  +1.35%  '),
    @('E30', '  -0.42%  '),
    @('E31', '  -0.18%  '),
    @('D32', '3.564.56'),
    @('E32', '  +0.15%  '),
    @('E33', '  +4.00%  '),
    @('D34', '24.79'),
    @('E34', '  -2.53%  '),
    @('E35', '  -3.01%  '),
    @('D37', '7.72'),
    @('E37', '  -1.25%  '),
    @('E38', '  -2.92%  '),
    @('E39', '  -5.15%  '),
    @('D40', '175.21'),
    @('E40', '  -0.01%  '),
    @('D41', '0.0843'),
    @('E41', '  -0.76%  '),
    @('D42', '5.17'),
    @('E42', '  -0.56%  '),
    @('D43', '0.875'),
    @('E43', '  -1.35%  '),
    @('D44', '45.74'),
    @('E44', '  -0.55%  '),
    @('E45', '  -4.37%  '),
    @('D46', '1.00'),
    @('E46', '  +0.00%  '),
    @('D47', '2.48'),
    @('E47', '  +0.31%  '),
    @('D48', '23.51'),
    @('E48', '  +0.80%  '),
    @('B49', 'ONDO'),
    @('C49', 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'),
    @('D49', '1.14'),
    @('E49', '  -4.76%  '),
    @('B50', 'InjectiveProtocol'),
    @('C50', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
    @('D50', '24.13'),
    @('E50', '  -3.09%  '),
    @('D51', '7.11'),
    @('E51', '  -0.21%  ')
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $ws.Range($ref).Value = "'" + $val
}
